# Update "想去人数" (interested-count) values in column F
# on sheets "展览" (Exhibitions) and "全部类型" (All Types).
$wb = $excel.ActiveWorkbook

# Row => New Value, for sheet "展览"
$sheet1Updates = @{
    4  = 14481
    5  = 17162
    7  = 141
    8  = 48
    17 = 13
    18 = 126
    20 = 1306
    25 = 7068
    26 = 976
    28 = 1155
    31 = 44
    32 = 5824
    36 = 4992
}

# Row => New Value, for sheet "全部类型"
$sheet4Updates = @{
    4  = 14481
    5  = 17162
    7  = 141
    8  = 48
    17 = 13
    18 = 126
    20 = 1306
    26 = 7068
    27 = 976
    29 = 1155
    32 = 44
    34 = 5824
    38 = 4992
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
